$d = $word.ActiveDocument

# Before:
#   Paragraph 1: "buy car"  (Harlow Solid Italic, 96pt)
#   Paragraph 2: <empty trailing paragraph>
#
# After (per the diff):
#   Paragraph 1: "buy car"  (unchanged)
#   Paragraph 2 (new): "Kms driven " styled with
#       rFonts ascii/hAnsi = "Bahnschrift SemiBold Condensed",
#       rFonts cs          = "Dreaming Outloud Pro",
#       sz/szCs = 36 (18pt)
#     applied to BOTH the paragraph mark (pPr/rPr) and the run (r/rPr).
#   Paragraph 3: <empty trailing paragraph>  (the original one, unchanged)

$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range

# Collapsed insertion point right after the "buy car" text but before its
# paragraph mark -- inserting a full <w:p> fragment here splits the story
# into a new paragraph positioned between paragraph 1 and the trailing
# empty paragraph, and (unlike the Font.* property setters) correctly
# stamps the new paragraph's own mark formatting (pPr/rPr) as well as the
# run formatting in one shot.
$insPoint = $d.Range($r1.End - 1, $r1.End - 1)

$newParaRPr = '<w:rFonts w:ascii="Bahnschrift SemiBold Condensed" w:hAnsi="Bahnschrift SemiBold Condensed" w:cs="Dreaming Outloud Pro"/><w:sz w:val="36"/><w:szCs w:val="36"/>'

$xmlFragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body>' + `
  '<w:p>' + `
  '<w:pPr><w:rPr>' + $newParaRPr + '</w:rPr></w:pPr>' + `
  '<w:r><w:rPr>' + $newParaRPr + '</w:rPr><w:t xml:space="preserve">Kms driven </w:t></w:r>' + `
  '</w:p>' + `
  '</w:body>' + `
  '</w:document>' + `
  '</pkg:xmlData></pkg:part></pkg:package>'

$insPoint.InsertXML($xmlFragment)
